$d = $word.ActiveDocument

# The sender address block at the top of the report currently reads:
#   {SenderName}
#   {SenderStreet}
#   {SenderCity}
#
# A department line was missing. Turn the {SenderStreet} line into the new
# {SenderDepartment} line, then insert a fresh {SenderStreet} line right
# after it, so the block becomes:
#   {SenderName}
#   {SenderDepartment}
#   {SenderStreet}
#   {SenderCity}

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("{SenderStreet}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $targetPara = $find.Parent.Paragraphs(1)
    $targetRange = $targetPara.Range
    $targetRange.Text = "{SenderDepartment}"
    $targetRange.InsertParagraphAfter()

    $newPara = $targetPara.Next()
    $newPara.Range.Text = "{SenderStreet}"
}
